# Edit script: corrects logic problems in the measurement data (column A recomputed,
# column C adjusted accordingly) and refreshes workbook window/path metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections (columns A and C) ---
$ws.Range("A19").Value = 12.02
$ws.Range("C19").Value = 134
$ws.Range("A21").Value = 26.740000000000002
$ws.Range("C21").Value = 131
$ws.Range("A22").Value = 40.72
$ws.Range("C22").Value = 138
$ws.Range("A29").Value = 7.04
$ws.Range("C29").Value = 141
$ws.Range("A30").Value = 4.55
$ws.Range("C30").Value = 140
$ws.Range("A32").Value = 17.96
$ws.Range("C32").Value = 131
$ws.Range("A35").Value = 9.5500000000000007
$ws.Range("C35").Value = 140
$ws.Range("A36").Value = 4.8
$ws.Range("C36").Value = 140
$ws.Range("A44").Value = 29.04
$ws.Range("C44").Value = 131
$ws.Range("A47").Value = 60.38
$ws.Range("C47").Value = 140
$ws.Range("A49").Value = 13.469999999999999
$ws.Range("C49").Value = 126
$ws.Range("A51").Value = 12.82
$ws.Range("C51").Value = 124
$ws.Range("A55").Value = 15.110000000000001
$ws.Range("C55").Value = 138
$ws.Range("A56").Value = 19.88
$ws.Range("C56").Value = 134
$ws.Range("A59").Value = 11.93
$ws.Range("C59").Value = 139
$ws.Range("A64").Value = 10.09
$ws.Range("C64").Value = 141
$ws.Range("A65").Value = 4.37
$ws.Range("C65").Value = 126
$ws.Range("A67").Value = 18.2
$ws.Range("C67").Value = 124
$ws.Range("A75").Value = 22.68
$ws.Range("C75").Value = 129
$ws.Range("A76").Value = 22.46
$ws.Range("C76").Value = 138
$ws.Range("A85").Value = 16.61
$ws.Range("C85").Value = 141
$ws.Range("A86").Value = 24.27
$ws.Range("C86").Value = 139
$ws.Range("A89").Value = 20.71
$ws.Range("C89").Value = 126
$ws.Range("A90").Value = 17.97
$ws.Range("C90").Value = 114
$ws.Range("A92").Value = 7.48
$ws.Range("C92").Value = 137
$ws.Range("A93").Value = 9.379999999999999
$ws.Range("C93").Value = 128
$ws.Range("A97").Value = 10.9
$ws.Range("C97").Value = 139
$ws.Range("A98").Value = 9.7199999999999989
$ws.Range("C98").Value = 139

# --- Workbook window metadata (maximize window to match the new layout) ---
$win = $excel.ActiveWindow
$win.WindowState = -4143
$win.Width = 25800
$win.Height = 13200
